$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was bumped by one day
# (45189 -> 45190) for every data row (rows 2 through 411).
$ws.Range("C2:C411").Value = 45190
